# Finish implementing the "Add Borrower" feature: append a new Borrower
# record (Barack Obama) to the Borrower table on Sheet1, turn the new
# emailAddress cell into a live mailto: hyperlink, and leave the
# selection where the data-entry form would have left it.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- New borrower row (row 6: bid, password, name, address, phone,
#     emailAddress, sintOrStdNo, expiryDate, type) -----------------------

$ws.Range("A6:G6").HorizontalAlignment = -4108   # xlCenter

$ws.Range("A6").Value = 1
$ws.Range("B6").Value = "password4"
$ws.Range("C6").Value = "Barack Obama"
$ws.Range("D6").Value = "Washington DC"
$ws.Range("E6").Value = 98765432

# emailAddress cell: add the mailto hyperlink, then set its display text.
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "mailto:obama@gmail.com")
$ws.Range("F6").HorizontalAlignment = -4108      # xlCenter
$ws.Range("F6").Value = "obama@gmail.com"

$ws.Range("G6").Value = 999999

# expiryDate: text-formatted like the other rows in the table.
$ws.Range("H6").NumberFormat = "@"
$ws.Range("H6").Value = "69-12-12"

# type
$ws.Range("I6").HorizontalAlignment = -4108      # xlCenter
$ws.Range("I6").Value = "Stu"

# --- Leave the selection where data entry finished -----------------------
$ws.Range("F8").Select()
